# The commit adds one new weekly price record for "Pomelo" (Start Ruby,
# Primera) dated 2021-09-08, inserted as a new row before the existing
# row 69, pushing the rest of the data table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 69 (shifts rows 69:129 down to 70:130)
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new record's data
$ws.Cells.Item(69, 1).Value = 10
$ws.Cells.Item(69, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(69, 3).Value = "La Araucanía"
$ws.Cells.Item(69, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(69, 5).Value = 9
$ws.Cells.Item(69, 6).Value = "Fruta"
$ws.Cells.Item(69, 7).Value = 100102
$ws.Cells.Item(69, 8).Value = "Cítricos"
$ws.Cells.Item(69, 9).Value = 100102006
$ws.Cells.Item(69, 10).Value = "Pomelo"
$ws.Cells.Item(69, 11).Value = "Start Ruby"
$ws.Cells.Item(69, 12).Value = "Primera"
$ws.Cells.Item(69, 13).Value = 80
$ws.Cells.Item(69, 14).Value = 12000
$ws.Cells.Item(69, 15).Value = 12000
$ws.Cells.Item(69, 16).Value = 12000
$ws.Cells.Item(69, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(69, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(69, 19).Value = 800
$ws.Cells.Item(69, 20).Value = 15
